$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 37450
$ws.Range("J3").Value = 37450
$ws.Range("L3").Value = 37450
$ws.Range("N3").Value = -37678
$ws.Range("H17").Value = 5102.769
$ws.Range("J17").Value = 6117.7144
$ws.Range("L17").Value = 18353.1432
$ws.Range("N17").Value = -18689.1432
$ws.Range("H32").Value = 8712.9
$ws.Range("I32").Value = 9655.333000000001
$ws.Range("K32").Value = 9655.333000000001
$ws.Range("M32").Value = -9329.333000000001
$ws.Range("H41").Value = 2886.25
$ws.Range("I41").Value = 2625.5557
$ws.Range("K41").Value = 2625.5557
$ws.Range("M41").Value = -2185.5557
$ws.Range("H62").Value = 9403.429
$ws.Range("J62").Value = 9636.333000000001
$ws.Range("L62").Value = 9636.333000000001
$ws.Range("N62").Value = -10884.333
$ws.Range("H65").Value = 9403.429
$ws.Range("J65").Value = 9636.333000000001
$ws.Range("L65").Value = 48181.665
$ws.Range("N65").Value = -54421.665
$ws.Range("H102").Value = 37450
$ws.Range("J102").Value = 37450
$ws.Range("L102").Value = 37450
$ws.Range("N102").Value = -43940
$ws.Range("H113").Value = 8072.273
$ws.Range("I113").Value = 2500
$ws.Range("K113").Value = 2500
$ws.Range("M113").Value = 754
$ws.Range("H129").Value = 1727.4286
$ws.Range("I129").Value = 932.55554
$ws.Range("J129").Value = 3158.2
$ws.Range("K129").Value = 2797.66662
$ws.Range("L129").Value = 9474.599999999999
$ws.Range("M129").Value = 2202.33338
$ws.Range("N129").Value = -19474.6
$ws.Range("H135").Value = 1292.4286
$ws.Range("I135").Value = 1733
$ws.Range("J135").Value = 705
$ws.Range("K135").Value = 15597
$ws.Range("L135").Value = 6345
$ws.Range("M135").Value = -13062
$ws.Range("N135").Value = -11415

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = $null
$ws.Range("H61").Value = 5196.515
$ws.Range("I61").Value = 5316.2
$ws.Range("K61").Value = 5316.2
$ws.Range("M61").Value = -5104.2
$ws.Range("H74").Value = 25643306
$ws.Range("I74").Value = 30305362
$ws.Range("K74").Value = 30305362
$ws.Range("M74").Value = -30304488
$ws.Range("H77").Value = 25643306
$ws.Range("I77").Value = 30305362
$ws.Range("K77").Value = 151526810
$ws.Range("M77").Value = -151522442
$ws.Range("H92").Value = 41998
$ws.Range("J92").Value = 41998
$ws.Range("L92").Value = 41998
$ws.Range("N92").Value = -46990
$ws.Range("H110").Value = 4167.294
$ws.Range("I110").Value = 1841.6923
$ws.Range("K110").Value = 1841.6923
$ws.Range("M110").Value = 203.3077000000001
$ws.Range("H135").Value = 69999
$ws.Range("J135").Value = 69999
$ws.Range("L135").Value = 69999
$ws.Range("N135").Value = -80139
$ws.Range("H136").Value = 5196.515
$ws.Range("I136").Value = 5316.2
$ws.Range("K136").Value = 15948.6
$ws.Range("M136").Value = -13398.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 666.4666999999999
$ws.Range("I94").Value = 642.6429000000001
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 642.6429000000001
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -191.6429000000001
$ws.Range("N94").Value = -1902
$ws.Range("H105").Value = 30500.1
$ws.Range("I105").Value = 35417.168
$ws.Range("K105").Value = 35417.168
$ws.Range("M105").Value = -33670.168
$ws.Range("H107").Value = 998.5
$ws.Range("I107").Value = 998.5
$ws.Range("K107").Value = 998.5
$ws.Range("M107").Value = 921.5
$ws.Range("H134").Value = 2702.6667
$ws.Range("I134").Value = 1248.25
$ws.Range("K134").Value = 3744.75
$ws.Range("M134").Value = -1209.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 62342.79
$ws.Range("I31").Value = 7720.625
$ws.Range("K31").Value = 7720.625
$ws.Range("M31").Value = -7425.625
$ws.Range("H34").Value = 62342.79
$ws.Range("I34").Value = 7720.625
$ws.Range("K34").Value = 7720.625
$ws.Range("M34").Value = -7518.625
$ws.Range("H58").Value = 4025.2273
$ws.Range("I58").Value = 1707.7059
$ws.Range("K58").Value = 1707.7059
$ws.Range("M58").Value = -1504.7059
$ws.Range("H134").Value = 5050.2
$ws.Range("I134").Value = 2738
$ws.Range("J134").Value = 7362.4
$ws.Range("K134").Value = 8214
$ws.Range("L134").Value = 22087.2
$ws.Range("M134").Value = -5679
$ws.Range("N134").Value = -27157.2
$ws.Range("H136").Value = 4025.2273
$ws.Range("I136").Value = 1707.7059
$ws.Range("K136").Value = 5123.1177
$ws.Range("M136").Value = -2573.1177

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2797.5293
$ws.Range("I34").Value = 1198
$ws.Range("J34").Value = 5730
$ws.Range("K34").Value = 3594
$ws.Range("L34").Value = 17190
$ws.Range("M34").Value = -3510
$ws.Range("N34").Value = -17358
$ws.Range("H38").Value = 36.1875
$ws.Range("I38").Value = 68.5
$ws.Range("J38").Value = 16.8
$ws.Range("K38").Value = 205.5
$ws.Range("L38").Value = 50.40000000000001
$ws.Range("M38").Value = 141.5
$ws.Range("N38").Value = -744.4
$ws.Range("H39").Value = 2983.2856
$ws.Range("J39").Value = 4245.75
$ws.Range("L39").Value = 12737.25
$ws.Range("N39").Value = -13325.25
$ws.Range("H55").Value = 6721
$ws.Range("I55").Value = 1582.2
$ws.Range("J55").Value = 9932.75
$ws.Range("K55").Value = 4746.6
$ws.Range("L55").Value = 29798.25
$ws.Range("M55").Value = -4569.6
$ws.Range("N55").Value = -30152.25
$ws.Range("H93").Value = 19027
$ws.Range("J93").Value = 19027
$ws.Range("L93").Value = 57081
$ws.Range("N93").Value = -60825
$ws.Range("H113").Value = 1325.7
$ws.Range("J113").Value = 1196
$ws.Range("L113").Value = 3588
$ws.Range("N113").Value = -7928

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 34500
$ws.Range("J12").Value = 4000
$ws.Range("L12").Value = 4000
$ws.Range("N12").Value = -4280
$ws.Range("H15").Value = 32000
$ws.Range("J15").Value = 32000
$ws.Range("L15").Value = 32000
$ws.Range("N15").Value = -32576
$ws.Range("H81").Value = 32000
$ws.Range("J81").Value = 32000
$ws.Range("L81").Value = 32000
$ws.Range("N81").Value = -33996
$ws.Range("H84").Value = 32000
$ws.Range("J84").Value = 32000
$ws.Range("L84").Value = 96000
$ws.Range("N84").Value = -105984
$ws.Range("H101").Value = 26116.666
$ws.Range("J101").Value = 26116.666
$ws.Range("L101").Value = 26116.666
$ws.Range("N101").Value = -32606.666
$ws.Range("H126").Value = 3280.3
$ws.Range("I126").Value = 2452.2104
$ws.Range("K126").Value = 7356.6312
$ws.Range("M126").Value = -4886.6312

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 27668.334
$ws.Range("I20").Value = 26502.5
$ws.Range("K20").Value = 26502.5
$ws.Range("M20").Value = -26276.5
$ws.Range("H100").Value = 12470.826
$ws.Range("J100").Value = 13710
$ws.Range("L100").Value = 13710
$ws.Range("N100").Value = -14792
$ws.Range("H116").Value = 259649
$ws.Range("J116").Value = 259649
$ws.Range("L116").Value = 259649
$ws.Range("N116").Value = -268827
$ws.Range("H132").Value = 6201.3
$ws.Range("I132").Value = 2752.1667
$ws.Range("J132").Value = 11375
$ws.Range("K132").Value = 8256.500100000001
$ws.Range("L132").Value = 34125
$ws.Range("M132").Value = -5726.500100000001
$ws.Range("N132").Value = -39185

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 465.08
$ws.Range("I107").Value = 390.66666
$ws.Range("K107").Value = 1171.99998
$ws.Range("M107").Value = 748.0000199999999
$ws.Range("H136").Value = 2633.7727
$ws.Range("I136").Value = 1854.1904
$ws.Range("K136").Value = 5562.5712
$ws.Range("M136").Value = -3012.5712
